$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 156, shifting existing rows 156:194 down to 157:195
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row 156 with the new weekly record
$ws.Range("A156").Value = 4
$ws.Range("B156").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C156").Value = "Los Lagos"
$ws.Range("D156").Value = 44551
$ws.Range("E156").Value = 10
$ws.Range("F156").Value = 100112003
$ws.Range("G156").Value = "Ajo"
$ws.Range("H156").Value = "Chino"
$ws.Range("I156").Value = "Primera"
$ws.Range("J156").Value = 240
$ws.Range("K156").Value = 21000
$ws.Range("L156").Value = 22000
$ws.Range("M156").Value = 21500
$ws.Range("N156").Value = "$/caja 10 kilos"
$ws.Range("O156").Value = "China"
$ws.Range("P156").Value = 2150
$ws.Range("Q156").Value = 10
$ws.Range("R156").Value = "Hortaliza"
